# Update scripts wuth new tpm
# The ligand (Gnas) average/total expression values for the "ECs" sending
# cluster were recomputed from the new TPM matrix. That change cascades
# into the derived specificity columns (I, J) for every row that shares
# the same ligand/receptor pair, and into the edge-weight columns
# (Q, R, S, T) for every row, since those are products of the ligand and
# receptor side values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> FAPs) ---
$ws.Range("G2").Value = 95.63567833333333
$ws.Range("H2").Value = 286.907035
$ws.Range("I2").Value = 0.2808828217467972
$ws.Range("J2").Value = 0.2808828217467972
$ws.Range("Q2").Value = 6.351898604983889
$ws.Range("R2").Value = 57.16708744485501
$ws.Range("S2").Value = 0.00481398502971861
$ws.Range("T2").Value = 0.004813985029718611

# --- Row 3 (ECs -> MuSCs) ---
$ws.Range("G3").Value = 95.63567833333333
$ws.Range("H3").Value = 286.907035
$ws.Range("I3").Value = 0.2808828217467972
$ws.Range("J3").Value = 0.2808828217467972
$ws.Range("Q3").Value = 364.2639617691617
$ws.Range("R3").Value = 3278.375655922455
$ws.Range("S3").Value = 0.2760688367170785
$ws.Range("T3").Value = 0.2760688367170786

# --- Row 4 (FAPs -> FAPs) ---
$ws.Range("I4").Value = 0.392628215788982
$ws.Range("J4").Value = 0.392628215788982
$ws.Range("R4").Value = 79.91023233724
$ws.Range("S4").Value = 0.006729163219376693
$ws.Range("T4").Value = 0.006729163219376693

# --- Row 5 (FAPs -> MuSCs) ---
$ws.Range("I5").Value = 0.392628215788982
$ws.Range("J5").Value = 0.392628215788982
$ws.Range("S5").Value = 0.3858990525696053
$ws.Range("T5").Value = 0.3858990525696053

# --- Row 6 (MuSCs -> FAPs) ---
$ws.Range("I6").Value = 0.3264889624642208
$ws.Range("J6").Value = 0.3264889624642208
$ws.Range("R6").Value = 66.44914399143001
$ws.Range("S6").Value = 0.005595618015714561
$ws.Range("T6").Value = 0.005595618015714562

# --- Row 7 (MuSCs -> MuSCs) ---
$ws.Range("I7").Value = 0.3264889624642208
$ws.Range("J7").Value = 0.3264889624642208
$ws.Range("S7").Value = 0.3208933444485062
$ws.Range("T7").Value = 0.3208933444485063
